$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.150.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.42%  "
$ws.Range("D3").Value = "'1.654.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'216.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.77%  "
$ws.Range("D6").Value = "'0.5124"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.96%  "
$ws.Range("D7").Value = "'1.012"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "'0.2596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").Value = "'0.06439"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").Value = "'19.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").Value = "'0.07813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "'1.664.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").Value = "'4.281"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.06%  "
$ws.Range("D14").Value = "'1.884.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").Value = "'0.5492"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.17%  "
$ws.Range("D16").Value = "'0.0₅7991"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").Value = "'63.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.55%  "
$ws.Range("D18").Value = "'26.186.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.27%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "'207.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.32%  "
$ws.Range("D21").Value = "'4.397"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.34%  "
$ws.Range("D22").Value = "'10.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.28%  "
$ws.Range("D23").Value = "'6.046"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'1.851"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.22%  "
$ws.Range("D26").Value = "'144.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").Value = "'0.1171"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D28").Value = "'6.954"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.73%  "
$ws.Range("D29").Value = "'15.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.31%  "
$ws.Range("D30").Value = "'0.05080"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.25%  "
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("D32").Value = "'3.341"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").Value = "'3.242"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("D34").Value = "'1.551"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.21%  "
$ws.Range("D35").Value = "'2.740"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.78%  "
$ws.Range("D36").Value = "'2.360"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").Value = "'0.9201"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("D38").Value = "'1.172.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").Value = "'0.5700"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("D40").Value = "'0.01585"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("D41").Value = "'1.012"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "'2.568"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "'5.650"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("D44").Value = "'0.8267"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "'100.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "'1.795.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.26%  "
$ws.Range("E47").Value = "  -4.23%  "
$ws.Range("D48").Value = "'0.4559"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "'1.010"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'55.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.65%  "
$ws.Range("D51").Value = "'7.859"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.67%  "
